# Auto-applied profit-column corrections (H:N) across multiple Garuda_Profits sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Cells.Item(34, 8).Value = 14200.5
$ws.Cells.Item(34, 9).Value = 1030.8
$ws.Cells.Item(34, 10).Value = 80049
$ws.Cells.Item(34, 11).Value = 1030.8
$ws.Cells.Item(34, 12).Value = 80049
$ws.Cells.Item(34, 13).Value = -827.8
$ws.Cells.Item(34, 14).Value = -80455

# Row 36
$ws.Cells.Item(36, 8).Value = 14200.5
$ws.Cells.Item(36, 9).Value = 1030.8
$ws.Cells.Item(36, 10).Value = 80049
$ws.Cells.Item(36, 11).Value = 1030.8
$ws.Cells.Item(36, 12).Value = 80049
$ws.Cells.Item(36, 13).Value = -315.8
$ws.Cells.Item(36, 14).Value = -81479

# Row 70
$ws.Cells.Item(70, 8).Value = 32601.188
$ws.Cells.Item(70, 9).Value = 84341.5
$ws.Cells.Item(70, 10).Value = 1557
$ws.Cells.Item(70, 11).Value = 253024.5
$ws.Cells.Item(70, 12).Value = 4671
$ws.Cells.Item(70, 13).Value = -252754.5
$ws.Cells.Item(70, 14).Value = -5211

# Row 73
$ws.Cells.Item(73, 8).Value = 32601.188
$ws.Cells.Item(73, 9).Value = 84341.5
$ws.Cells.Item(73, 10).Value = 1557
$ws.Cells.Item(73, 11).Value = 253024.5
$ws.Cells.Item(73, 12).Value = 4671
$ws.Cells.Item(73, 13).Value = -252088.5
$ws.Cells.Item(73, 14).Value = -6543

# Row 116
$ws.Cells.Item(116, 8).Value = 1960.3572
$ws.Cells.Item(116, 9).Value = 1799
$ws.Cells.Item(116, 10).Value = 1972.7693
$ws.Cells.Item(116, 11).Value = 1799
$ws.Cells.Item(116, 12).Value = 1972.7693
$ws.Cells.Item(116, 13).Value = 1643
$ws.Cells.Item(116, 14).Value = -8856.7693

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Cells.Item(43, 8).Value = 5500
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 5500
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 5500
$ws.Cells.Item(43, 14).Value = -6126

# Row 61
$ws.Cells.Item(61, 8).Value = 1662.2927
$ws.Cells.Item(61, 9).Value = 912.06665
$ws.Cells.Item(61, 10).Value = 3708.3635
$ws.Cells.Item(61, 11).Value = 912.06665
$ws.Cells.Item(61, 12).Value = 3708.3635
$ws.Cells.Item(61, 13).Value = -700.06665
$ws.Cells.Item(61, 14).Value = -4132.363499999999

# Row 74
$ws.Cells.Item(74, 8).Value = 653.875
$ws.Cells.Item(74, 9).Value = 499.8
$ws.Cells.Item(74, 10).Value = 910.6667
$ws.Cells.Item(74, 11).Value = 499.8
$ws.Cells.Item(74, 12).Value = 910.6667
$ws.Cells.Item(74, 13).Value = 374.2
$ws.Cells.Item(74, 14).Value = -2658.6667

# Row 77
$ws.Cells.Item(77, 8).Value = 653.875
$ws.Cells.Item(77, 9).Value = 499.8
$ws.Cells.Item(77, 10).Value = 910.6667
$ws.Cells.Item(77, 11).Value = 2499
$ws.Cells.Item(77, 12).Value = 4553.3335
$ws.Cells.Item(77, 13).Value = 1869
$ws.Cells.Item(77, 14).Value = -13289.3335

# Row 136
$ws.Cells.Item(136, 8).Value = 1662.2927
$ws.Cells.Item(136, 9).Value = 912.06665
$ws.Cells.Item(136, 10).Value = 3708.3635
$ws.Cells.Item(136, 11).Value = 2736.19995
$ws.Cells.Item(136, 12).Value = 11125.0905
$ws.Cells.Item(136, 13).Value = -186.1999500000002
$ws.Cells.Item(136, 14).Value = -16225.0905

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Cells.Item(99, 8).Value = 2893.6924
$ws.Cells.Item(99, 9).Value = 3200
$ws.Cells.Item(99, 10).Value = 2757.5557
$ws.Cells.Item(99, 11).Value = 3200
$ws.Cells.Item(99, 12).Value = 2757.5557
$ws.Cells.Item(99, 13).Value = -1702
$ws.Cells.Item(99, 14).Value = -5753.5557

# Row 126
$ws.Cells.Item(126, 8).Value = 2893.6924
$ws.Cells.Item(126, 9).Value = 3200
$ws.Cells.Item(126, 10).Value = 2757.5557
$ws.Cells.Item(126, 11).Value = 9600
$ws.Cells.Item(126, 12).Value = 8272.667099999999
$ws.Cells.Item(126, 13).Value = -7130
$ws.Cells.Item(126, 14).Value = -13212.6671

# Row 132
$ws.Cells.Item(132, 8).Value = 3678709.5
$ws.Cells.Item(132, 9).Value = 2481.8125
$ws.Cells.Item(132, 10).Value = 6946467.5
$ws.Cells.Item(132, 11).Value = 7445.4375
$ws.Cells.Item(132, 12).Value = 20839402.5
$ws.Cells.Item(132, 13).Value = -4915.4375
$ws.Cells.Item(132, 14).Value = -20844462.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Cells.Item(8, 8).Value = 172
$ws.Cells.Item(8, 9).Value = 172
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 516
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = -377

# Row 86
$ws.Cells.Item(86, 8).Value = 619.9167
$ws.Cells.Item(86, 9).Value = 689.8333
$ws.Cells.Item(86, 10).Value = 550
$ws.Cells.Item(86, 11).Value = 2069.4999
$ws.Cells.Item(86, 12).Value = 1650
$ws.Cells.Item(86, 13).Value = -883.4998999999998
$ws.Cells.Item(86, 14).Value = -4022

# Row 89
$ws.Cells.Item(89, 8).Value = 619.9167
$ws.Cells.Item(89, 9).Value = 689.8333
$ws.Cells.Item(89, 10).Value = 550
$ws.Cells.Item(89, 11).Value = 6208.4997
$ws.Cells.Item(89, 12).Value = 4950
$ws.Cells.Item(89, 13).Value = -280.4997000000003
$ws.Cells.Item(89, 14).Value = -16806

# Row 92
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).ClearContents()
$ws.Cells.Item(92, 14).ClearContents()

# Row 122
$ws.Cells.Item(122, 8).Value = 1002.2105
$ws.Cells.Item(122, 9).Value = 396.42856
$ws.Cells.Item(122, 10).Value = 1355.5834
$ws.Cells.Item(122, 11).Value = 3567.85704
$ws.Cells.Item(122, 12).Value = 12200.2506
$ws.Cells.Item(122, 13).Value = -1117.85704
$ws.Cells.Item(122, 14).Value = -17100.2506

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 449.22726
$ws.Cells.Item(107, 9).Value = 483.33334
$ws.Cells.Item(107, 10).Value = 376.14285
$ws.Cells.Item(107, 11).Value = 483.33334
$ws.Cells.Item(107, 12).Value = 376.14285
$ws.Cells.Item(107, 13).Value = 1436.66666
$ws.Cells.Item(107, 14).Value = -4216.14285

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 1860
$ws.Cells.Item(68, 9).Value = 1280
$ws.Cells.Item(68, 10).Value = 2150
$ws.Cells.Item(68, 11).Value = 1280
$ws.Cells.Item(68, 12).Value = 2150
$ws.Cells.Item(68, 13).Value = -531
$ws.Cells.Item(68, 14).Value = -3648

# Row 71
$ws.Cells.Item(71, 8).Value = 1860
$ws.Cells.Item(71, 9).Value = 1280
$ws.Cells.Item(71, 10).Value = 2150
$ws.Cells.Item(71, 11).Value = 6400
$ws.Cells.Item(71, 12).Value = 10750
$ws.Cells.Item(71, 13).Value = -2656
$ws.Cells.Item(71, 14).Value = -18238

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Cells.Item(39, 8).Value = 2229.3333
$ws.Cells.Item(39, 9).Value = 894
$ws.Cells.Item(39, 10).Value = 4900
$ws.Cells.Item(39, 11).Value = 894
$ws.Cells.Item(39, 12).Value = 4900
$ws.Cells.Item(39, 13).Value = -481
$ws.Cells.Item(39, 14).Value = -5726

# Row 42
$ws.Cells.Item(42, 8).Value = 80049
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 80049
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 80049
$ws.Cells.Item(42, 14).Value = -80805

# Row 43
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).ClearContents()

# Row 62
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).ClearContents()

# Row 65
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).ClearContents()

# Row 107
$ws.Cells.Item(107, 8).Value = 490.3
$ws.Cells.Item(107, 9).Value = 523.8570999999999
$ws.Cells.Item(107, 10).Value = 412
$ws.Cells.Item(107, 11).Value = 1571.5713
$ws.Cells.Item(107, 12).Value = 1236
$ws.Cells.Item(107, 13).Value = 348.4287000000002
$ws.Cells.Item(107, 14).Value = -5076

# Row 132
$ws.Cells.Item(132, 8).Value = 1399.2972
$ws.Cells.Item(132, 9).Value = 1055.5
$ws.Cells.Item(132, 10).Value = 3599.6
$ws.Cells.Item(132, 11).Value = 3166.5
$ws.Cells.Item(132, 12).Value = 10798.8
$ws.Cells.Item(132, 13).Value = -636.5
$ws.Cells.Item(132, 14).Value = -15858.8

# Row 136
$ws.Cells.Item(136, 8).Value = 3011.0356
$ws.Cells.Item(136, 9).Value = 3047.8462
$ws.Cells.Item(136, 10).Value = 2532.5
$ws.Cells.Item(136, 11).Value = 9143.5386
$ws.Cells.Item(136, 12).Value = 7597.5
$ws.Cells.Item(136, 13).Value = -6593.5386
$ws.Cells.Item(136, 14).Value = -12697.5
